$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 443, shifting the existing
# rows 443-452 down to become rows 446-455.
$ws.Rows.Item(443).Resize(3).Insert()

# New weekly price records (Hortaliza / Betarraga) to fill the 3 rows
# that were just inserted at 443-445.
$newRows = @(
    @{Row=443; D=44448; I="Primera"; J=52000; K=110; L=120; M=114; P=114},
    @{Row=444; D=44448; I="Segunda"; J=42000; K=90;  L=95;  M=92;  P=92},
    @{Row=445; D=44448; I="Tercera"; J=12000; K=60;  L=60;  M=60;  P=60}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row,1).Value  = 6
    $ws.Cells.Item($row,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($row,3).Value  = "Metropolitana"
    $ws.Cells.Item($row,4).Value  = $r.D
    $ws.Cells.Item($row,5).Value  = 13
    $ws.Cells.Item($row,6).Value  = 100114014
    $ws.Cells.Item($row,7).Value  = "Betarraga"
    $ws.Cells.Item($row,8).Value  = "Sin especificar"
    $ws.Cells.Item($row,9).Value  = $r.I
    $ws.Cells.Item($row,10).Value = $r.J
    $ws.Cells.Item($row,11).Value = $r.K
    $ws.Cells.Item($row,12).Value = $r.L
    $ws.Cells.Item($row,13).Value = $r.M
    $ws.Cells.Item($row,14).Value = "$/unidad"
    $ws.Cells.Item($row,15).Value = "Región Metropolitana"
    $ws.Cells.Item($row,16).Value = $r.P
    $ws.Cells.Item($row,17).Value = 1
    $ws.Cells.Item($row,18).Value = "Hortaliza"
}
